$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Hourly rows (HOUR 1-18, sheet rows 2-19): ACTUAL_ENERGY (B), TOTAL_BCQ_NOMINATION (C), WESM_EXPOSURE (D)
$ws.Range("B2").Value = 30928.199
$ws.Range("C2").Value = 66500
$ws.Range("D2").Value = -35571.801

$ws.Range("B3").Value = 29619.8635
$ws.Range("C3").Value = 63500
$ws.Range("D3").Value = -33880.1365

$ws.Range("B4").Value = 27793.1675
$ws.Range("C4").Value = 42500
$ws.Range("D4").Value = -14706.8325

$ws.Range("B5").Value = 26218.02
$ws.Range("C5").Value = 42500
$ws.Range("D5").Value = -16281.98

$ws.Range("B6").Value = 25305.1185
$ws.Range("C6").Value = 22500
$ws.Range("D6").Value = 2805.1185

$ws.Range("B7").Value = 25440.605
$ws.Range("C7").Value = 22500
$ws.Range("D7").Value = 2940.605

$ws.Range("B8").Value = 25564.52277227723
$ws.Range("C8").Value = 22500
$ws.Range("D8").Value = 3064.522772277229

$ws.Range("B9").Value = 26075.5535
$ws.Range("C9").Value = 22500
$ws.Range("D9").Value = 3575.553500000002

$ws.Range("B10").Value = 28306.4105
$ws.Range("C10").Value = 22500
$ws.Range("D10").Value = 5806.410500000002

$ws.Range("B11").Value = 26016.9995
$ws.Range("C11").Value = 32500
$ws.Range("D11").Value = -6483.000499999998

$ws.Range("B12").Value = 24997.2015
$ws.Range("C12").Value = 32500
$ws.Range("D12").Value = -7502.798500000001

$ws.Range("B13").Value = 24407.995
$ws.Range("C13").Value = 32500
$ws.Range("D13").Value = -8092.005000000001

$ws.Range("B14").Value = 23983.6615
$ws.Range("C14").Value = 32500
$ws.Range("D14").Value = -8516.338500000002

$ws.Range("B15").Value = 25520.364
$ws.Range("C15").Value = 52500
$ws.Range("D15").Value = -26979.636

$ws.Range("B16").Value = 26075.793
$ws.Range("C16").Value = 75000
$ws.Range("D16").Value = -48924.20699999999

$ws.Range("B17").Value = 25427.0925
$ws.Range("C17").Value = 75000
$ws.Range("D17").Value = -49572.9075

$ws.Range("B18").Value = 23201.19782374277
$ws.Range("C18").Value = 75000
$ws.Range("D18").Value = -51798.80217625723

$ws.Range("B19").Value = 20688.49828448191
$ws.Range("C19").Value = 75000
$ws.Range("D19").Value = -54311.50171551809

# Rows 24-25 only have TOTAL_BCQ_NOMINATION changed (no WESM_EXPOSURE cell present)
$ws.Range("C24").Value = 72000
$ws.Range("C25").Value = 52500
